$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44188
$ws.Range("O2").Value = "Región Metropolitana"

# Row 3
$ws.Range("D3").Value = 44160
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 9500
$ws.Range("P3").Value = 380

# Row 4
$ws.Range("D4").Value = 44162
$ws.Range("K4").Value = 7500
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7750
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 310

# Row 5
$ws.Range("D5").Value = 44384
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 12500
$ws.Range("O5").Value = "Región de Coquimbo"
$ws.Range("P5").Value = 500

# Row 6
$ws.Range("D6").Value = 44335
$ws.Range("O6").Value = "Provincia de Limarí"

# Row 8
$ws.Range("D8").Value = 44351
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("P8").Value = 620
